# Update row 8 (year 2025) metrics in the recorrencia anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1235
$ws.Range("D8").Value = 201
$ws.Range("E8").Value = 1034
$ws.Range("F8").Value = 8.244462674323216
$ws.Range("G8").Value = 83.7246963562753
$ws.Range("H8").Value = 16.27530364372469
